$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting (style) from the existing A90 cell down to A91:A94
$ws.Range("A90").Copy()
$ws.Range("A91:A94").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in dates and hour counts first (no shared-string impact)
$ws.Range("A91").Value = 42939
$ws.Range("C91").Value = 6

$ws.Range("A92").Value = 42940
$ws.Range("C92").Value = 2

$ws.Range("A93").Value = 42941
$ws.Range("C93").Value = 3

$ws.Range("A94").Value = 42942
$ws.Range("C94").Value = 3

# Now set the new shared strings in the same order they were introduced
# upstream: Monitor, Opravy chyb, revize, Backup, Vkladani prubehu signalu
$ws.Range("B90").Value = "Monitor"
$ws.Range("B92").Value = "Opravy chyb, revize"
$ws.Range("B94").Value = "Backup"
$ws.Range("B91").Value = "Vkládání průběhů signálu"

# Update the view so the new last row is visible/selected, matching the diff
[void]$ws.Range("A94:C94").Select()
$ws.Application.ActiveWindow.ScrollRow = 76
